# Auto-generated edit script: updates crypto price/volume table
# to reflect the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.451.76"
$ws.Range("E2").Value = "  -0.58%  "

# Row 3
$ws.Range("D3").Value = "3.779.85"
$ws.Range("E3").Value = "  +0.59%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "616.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "

# Row 7
$ws.Range("D7").Value = "3.779.20"
$ws.Range("E7").Value = "  +0.64%  "

# Row 8
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.96%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.167"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.35%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.52"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.14%  "

# Row 12
$ws.Range("E12").Value = "  -0.19%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.31%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000255"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.73%  "

# Row 15
$ws.Range("D15").Value = "4.405.53"
$ws.Range("E15").Value = "  +0.66%  "

# Row 16
$ws.Range("D16").Value = "3.779.48"
$ws.Range("E16").Value = "  +0.72%  "

# Row 17
$ws.Range("D17").Value = "69.497.60"
$ws.Range("E17").Value = "  -0.50%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.45%  "

# Row 19
$ws.Range("E19").Value = "  -3.29%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "509.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.07%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.67%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.83%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.733"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.52%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.10%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.55%  "

# Row 26
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000144"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.52%  "

# Row 27
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.64%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.60%  "

# Row 29
$ws.Range("E29").Value = "  -0.01%  "

# Row 30
$ws.Range("E30").Value = "  +1.79%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.83%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.07%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.23%  "

# Row 34
$ws.Range("E34").Value = "  +0.94%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "

# Row 36
$ws.Range("E36").Value = "  -1.41%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.13%  "

# Row 38
$ws.Range("E38").Value = "  +7.07%  "

# Row 39
$ws.Range("E39").Value = "  +1.81%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "466.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.20%  "

# Row 41
$ws.Range("E41").Value = "  -1.83%  "

# Row 42
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.17%  "

# Row 43
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "49.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.25%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.96%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.81%  "

# Row 46
$ws.Range("D46").Value = "2.956.26"
$ws.Range("E46").Value = "  -1.83%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0363"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.60%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.85%  "

# Row 49
$ws.Range("E49").Value = "  +0.02%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "138.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.44%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.04%  "

